$d = $word.ActiveDocument
$d.Content.Find.Execute("runs successfully", $true, $false, $false, $false, $false, $true, 1, $false, "runs successfully.", 2)
